$d = $word.ActiveDocument

# Locate the Bibliografia entries paragraph (currently one run with all
# references concatenated with no separation) and replace its content with
# a version that inserts manual line breaks between / within entries.
$target = $null
foreach ($p in $d.Paragraphs) {
  if ($p.Range.Text -like "*IEZZI, G.*") {
    $target = $p.Range
  }
}

if ($target -eq $null) {
  throw "Could not locate the bibliography paragraph"
}

# Exclude the trailing paragraph-mark character so InsertXML replaces only
# the paragraph's content, not the paragraph itself.
$contentRange = $d.Range($target.Start, $target.End - 1)

$xmlFrag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">IEZZI, G.; MURAKAMI, C. Fundamentos de matemática elementar - Volume 1: Conjuntos e funções. São Paulo: Saraiva Didáticos, 2019. </w:t><w:br/><w:br/><w:t>CASTRUCCI, B.; GIOVANNI, J.R.; GIOVANNI JR., J.R. A Conquista da Matemática - 6º ano. São Paulo: FTD Educação, 2022.</w:t><w:br/><w:br/><w:t>FILHO, B. B.; SILVA, C. X. Matemática aula por aula. São Paulo: FTD, 2000.</w:t><w:br/><w:br/><w:t xml:space="preserve">WAGNER, E. Matemática 1. 1ª ed. Rio de Janeiro: FGV, 2011. </w:t><w:br/><w:br/><w:t>ADAMI, A. M.; DORNELLES FILHO, A. A.; LORANDI, M. M.: Pré-Cálculo. São Paulo: Bookman</w:t><w:br/><w:t>Editora, 2015</w:t><w:br/><w:br/><w:t>STEWART, J., REDLIN, L. e WATSON, S. Precalculus: Mathematics for Calculus.</w:t><w:br/><w:t>São Paulo: Cengage Learning, 7a ed., 2014.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$contentRange.InsertXML($xmlFrag)
